$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.04%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.62%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.047'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.40%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07963'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.73%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.867'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-4.59%'

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.764'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.11%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9200'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.20%'

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1338'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.08%'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1890'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.11%'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09092'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.05%'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03437'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.53%'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09850'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.30%'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001408'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-2.07%'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006062'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.68%'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.743'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.38%'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.117'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.82%'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '13.21%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3444'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.04%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.23%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.155'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.10%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2348'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-6.32%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04409'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.66%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001224'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.14%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004628'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.31%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001298'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.46%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004438'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.08%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01942'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.11%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05376'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '8.93%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007611'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.72%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01011'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.40%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1354'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.88%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002136'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.54%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01017'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-9.72%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006153'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.56%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.26%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.22'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.85%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001657'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '38.73%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.26%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.26%'
